$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(42,1).Value = "'2026-02-04"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).Value = "14:04:15"
$ws.Cells.Item(42,3).Value = "14:00"
$ws.Cells.Item(42,4).Value = "Bathroom"
$ws.Cells.Item(42,5).Value = "No Motion"
$ws.Cells.Item(42,6).Value = "Inactive"
$ws.Cells.Item(43,1).Value = "'2026-02-04"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).Value = "14:04:17"
$ws.Cells.Item(43,3).Value = "14:00"
$ws.Cells.Item(43,4).Value = "Bathroom"
$ws.Cells.Item(43,5).Value = "No Motion"
$ws.Cells.Item(43,6).Value = "Inactive"
$ws.Cells.Item(44,1).Value = "'2026-02-04"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).Value = "14:04:22"
$ws.Cells.Item(44,3).Value = "14:00"
$ws.Cells.Item(44,4).Value = "Bathroom"
$ws.Cells.Item(44,5).Value = "No Motion"
$ws.Cells.Item(44,6).Value = "Inactive"
$ws.Cells.Item(45,1).Value = "'2026-02-04"
$ws.Cells.Item(45,1).Style = "Normal"
$ws.Cells.Item(45,2).Value = "14:04:27"
$ws.Cells.Item(45,3).Value = "14:00"
$ws.Cells.Item(45,4).Value = "Bathroom"
$ws.Cells.Item(45,5).Value = "Motion Detected"
$ws.Cells.Item(45,6).Value = "Active"
$ws.Cells.Item(46,1).Value = "'2026-02-04"
$ws.Cells.Item(46,1).Style = "Normal"
$ws.Cells.Item(46,2).Value = "14:04:35"
$ws.Cells.Item(46,3).Value = "14:00"
$ws.Cells.Item(46,4).Value = "Bathroom"
$ws.Cells.Item(46,5).Value = "No Motion"
$ws.Cells.Item(46,6).Value = "Inactive"
$ws.Cells.Item(47,1).Value = "'2026-02-04"
$ws.Cells.Item(47,1).Style = "Normal"
$ws.Cells.Item(47,2).Value = "14:04:40"
$ws.Cells.Item(47,3).Value = "14:00"
$ws.Cells.Item(47,4).Value = "Bathroom"
$ws.Cells.Item(47,5).Value = "No Motion"
$ws.Cells.Item(47,6).Value = "Inactive"
$ws.Cells.Item(48,1).Value = "'2026-02-04"
$ws.Cells.Item(48,1).Style = "Normal"
$ws.Cells.Item(48,2).Value = "14:04:44"
$ws.Cells.Item(48,3).Value = "14:00"
$ws.Cells.Item(48,4).Value = "Bathroom"
$ws.Cells.Item(48,5).Value = "Motion Detected"
$ws.Cells.Item(48,6).Value = "Active"
$ws.Cells.Item(49,1).Value = "'2026-02-04"
$ws.Cells.Item(49,1).Style = "Normal"
$ws.Cells.Item(49,2).Value = "14:04:52"
$ws.Cells.Item(49,3).Value = "14:00"
$ws.Cells.Item(49,4).Value = "Bathroom"
$ws.Cells.Item(49,5).Value = "No Motion"
$ws.Cells.Item(49,6).Value = "Inactive"
$ws.Cells.Item(50,1).Value = "'2026-02-04"
$ws.Cells.Item(50,1).Style = "Normal"
$ws.Cells.Item(50,2).Value = "14:04:57"
$ws.Cells.Item(50,3).Value = "14:00"
$ws.Cells.Item(50,4).Value = "Bathroom"
$ws.Cells.Item(50,5).Value = "No Motion"
$ws.Cells.Item(50,6).Value = "Inactive"
$ws.Cells.Item(51,1).Value = "'2026-02-04"
$ws.Cells.Item(51,1).Style = "Normal"
$ws.Cells.Item(51,2).Value = "14:05:02"
$ws.Cells.Item(51,3).Value = "14:00"
$ws.Cells.Item(51,4).Value = "Bathroom"
$ws.Cells.Item(51,5).Value = "No Motion"
$ws.Cells.Item(51,6).Value = "Inactive"
$ws.Cells.Item(52,1).Value = "'2026-02-04"
$ws.Cells.Item(52,1).Style = "Normal"
$ws.Cells.Item(52,2).Value = "14:05:07"
$ws.Cells.Item(52,3).Value = "14:00"
$ws.Cells.Item(52,4).Value = "Bathroom"
$ws.Cells.Item(52,5).Value = "Motion Detected"
$ws.Cells.Item(52,6).Value = "Active"
$ws.Cells.Item(53,1).Value = "'2026-02-04"
$ws.Cells.Item(53,1).Style = "Normal"
$ws.Cells.Item(53,2).Value = "14:05:14"
$ws.Cells.Item(53,3).Value = "14:00"
$ws.Cells.Item(53,4).Value = "Bathroom"
$ws.Cells.Item(53,5).Value = "No Motion"
$ws.Cells.Item(53,6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(35,1).Value = "'2026-02-04"
$ws.Cells.Item(35,1).Style = "Normal"
$ws.Cells.Item(35,2).Value = "14:04:16"
$ws.Cells.Item(35,3).Value = "14:00"
$ws.Cells.Item(35,4).Value = "Bathroom"
$ws.Cells.Item(35,5).Value = "'76.6%"
$ws.Cells.Item(35,5).Style = "Normal"
$ws.Cells.Item(35,6).Value = "Active"
$ws.Cells.Item(36,1).Value = "'2026-02-04"
$ws.Cells.Item(36,1).Style = "Normal"
$ws.Cells.Item(36,2).Value = "14:04:18"
$ws.Cells.Item(36,3).Value = "14:00"
$ws.Cells.Item(36,4).Value = "Bathroom"
$ws.Cells.Item(36,5).Value = "'77.7%"
$ws.Cells.Item(36,5).Style = "Normal"
$ws.Cells.Item(36,6).Value = "Active"
$ws.Cells.Item(37,1).Value = "'2026-02-04"
$ws.Cells.Item(37,1).Style = "Normal"
$ws.Cells.Item(37,2).Value = "14:04:33"
$ws.Cells.Item(37,3).Value = "14:00"
$ws.Cells.Item(37,4).Value = "Bathroom"
$ws.Cells.Item(37,5).Value = "'76.7%"
$ws.Cells.Item(37,5).Style = "Normal"
$ws.Cells.Item(37,6).Value = "Active"
$ws.Cells.Item(38,1).Value = "'2026-02-04"
$ws.Cells.Item(38,1).Style = "Normal"
$ws.Cells.Item(38,2).Value = "14:04:38"
$ws.Cells.Item(38,3).Value = "14:00"
$ws.Cells.Item(38,4).Value = "Bathroom"
$ws.Cells.Item(38,5).Value = "'77.7%"
$ws.Cells.Item(38,5).Style = "Normal"
$ws.Cells.Item(38,6).Value = "Active"
$ws.Cells.Item(39,1).Value = "'2026-02-04"
$ws.Cells.Item(39,1).Style = "Normal"
$ws.Cells.Item(39,2).Value = "14:04:48"
$ws.Cells.Item(39,3).Value = "14:00"
$ws.Cells.Item(39,4).Value = "Bathroom"
$ws.Cells.Item(39,5).Value = "'78.2%"
$ws.Cells.Item(39,5).Style = "Normal"
$ws.Cells.Item(39,6).Value = "Active"
$ws.Cells.Item(40,1).Value = "'2026-02-04"
$ws.Cells.Item(40,1).Style = "Normal"
$ws.Cells.Item(40,2).Value = "14:04:53"
$ws.Cells.Item(40,3).Value = "14:00"
$ws.Cells.Item(40,4).Value = "Bathroom"
$ws.Cells.Item(40,5).Value = "'76.7%"
$ws.Cells.Item(40,5).Style = "Normal"
$ws.Cells.Item(40,6).Value = "Active"
$ws.Cells.Item(41,1).Value = "'2026-02-04"
$ws.Cells.Item(41,1).Style = "Normal"
$ws.Cells.Item(41,2).Value = "14:04:58"
$ws.Cells.Item(41,3).Value = "14:00"
$ws.Cells.Item(41,4).Value = "Bathroom"
$ws.Cells.Item(41,5).Value = "'77.6%"
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(41,6).Value = "Active"
$ws.Cells.Item(42,1).Value = "'2026-02-04"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).Value = "14:05:03"
$ws.Cells.Item(42,3).Value = "14:00"
$ws.Cells.Item(42,4).Value = "Bathroom"
$ws.Cells.Item(42,5).Value = "'76.7%"
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(42,6).Value = "Active"
$ws.Cells.Item(43,1).Value = "'2026-02-04"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).Value = "14:05:08"
$ws.Cells.Item(43,3).Value = "14:00"
$ws.Cells.Item(43,4).Value = "Bathroom"
$ws.Cells.Item(43,5).Value = "'77.6%"
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(43,6).Value = "Active"
$ws.Cells.Item(44,1).Value = "'2026-02-04"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).Value = "14:05:13"
$ws.Cells.Item(44,3).Value = "14:00"
$ws.Cells.Item(44,4).Value = "Bathroom"
$ws.Cells.Item(44,5).Value = "'76.4%"
$ws.Cells.Item(44,5).Style = "Normal"
$ws.Cells.Item(44,6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(35,1).Value = "'2026-02-04"
$ws.Cells.Item(35,1).Style = "Normal"
$ws.Cells.Item(35,2).Value = "14:04:16"
$ws.Cells.Item(35,3).Value = "14:00"
$ws.Cells.Item(35,4).Value = "Bathroom"
$ws.Cells.Item(35,5).Value = "24.7C"
$ws.Cells.Item(35,6).Value = "Active"
$ws.Cells.Item(36,1).Value = "'2026-02-04"
$ws.Cells.Item(36,1).Style = "Normal"
$ws.Cells.Item(36,2).Value = "14:04:18"
$ws.Cells.Item(36,3).Value = "14:00"
$ws.Cells.Item(36,4).Value = "Bathroom"
$ws.Cells.Item(36,5).Value = "24.8C"
$ws.Cells.Item(36,6).Value = "Active"
$ws.Cells.Item(37,1).Value = "'2026-02-04"
$ws.Cells.Item(37,1).Style = "Normal"
$ws.Cells.Item(37,2).Value = "14:04:33"
$ws.Cells.Item(37,3).Value = "14:00"
$ws.Cells.Item(37,4).Value = "Bathroom"
$ws.Cells.Item(37,5).Value = "24.8C"
$ws.Cells.Item(37,6).Value = "Active"
$ws.Cells.Item(38,1).Value = "'2026-02-04"
$ws.Cells.Item(38,1).Style = "Normal"
$ws.Cells.Item(38,2).Value = "14:04:38"
$ws.Cells.Item(38,3).Value = "14:00"
$ws.Cells.Item(38,4).Value = "Bathroom"
$ws.Cells.Item(38,5).Value = "24.8C"
$ws.Cells.Item(38,6).Value = "Active"
$ws.Cells.Item(39,1).Value = "'2026-02-04"
$ws.Cells.Item(39,1).Style = "Normal"
$ws.Cells.Item(39,2).Value = "14:04:48"
$ws.Cells.Item(39,3).Value = "14:00"
$ws.Cells.Item(39,4).Value = "Bathroom"
$ws.Cells.Item(39,5).Value = "24.8C"
$ws.Cells.Item(39,6).Value = "Active"
$ws.Cells.Item(40,1).Value = "'2026-02-04"
$ws.Cells.Item(40,1).Style = "Normal"
$ws.Cells.Item(40,2).Value = "14:04:53"
$ws.Cells.Item(40,3).Value = "14:00"
$ws.Cells.Item(40,4).Value = "Bathroom"
$ws.Cells.Item(40,5).Value = "24.8C"
$ws.Cells.Item(40,6).Value = "Active"
$ws.Cells.Item(41,1).Value = "'2026-02-04"
$ws.Cells.Item(41,1).Style = "Normal"
$ws.Cells.Item(41,2).Value = "14:04:58"
$ws.Cells.Item(41,3).Value = "14:00"
$ws.Cells.Item(41,4).Value = "Bathroom"
$ws.Cells.Item(41,5).Value = "24.8C"
$ws.Cells.Item(41,6).Value = "Active"
$ws.Cells.Item(42,1).Value = "'2026-02-04"
$ws.Cells.Item(42,1).Style = "Normal"
$ws.Cells.Item(42,2).Value = "14:05:03"
$ws.Cells.Item(42,3).Value = "14:00"
$ws.Cells.Item(42,4).Value = "Bathroom"
$ws.Cells.Item(42,5).Value = "24.8C"
$ws.Cells.Item(42,6).Value = "Active"
$ws.Cells.Item(43,1).Value = "'2026-02-04"
$ws.Cells.Item(43,1).Style = "Normal"
$ws.Cells.Item(43,2).Value = "14:05:08"
$ws.Cells.Item(43,3).Value = "14:00"
$ws.Cells.Item(43,4).Value = "Bathroom"
$ws.Cells.Item(43,5).Value = "24.8C"
$ws.Cells.Item(43,6).Value = "Active"
$ws.Cells.Item(44,1).Value = "'2026-02-04"
$ws.Cells.Item(44,1).Style = "Normal"
$ws.Cells.Item(44,2).Value = "14:05:13"
$ws.Cells.Item(44,3).Value = "14:00"
$ws.Cells.Item(44,4).Value = "Bathroom"
$ws.Cells.Item(44,5).Value = "24.8C"
$ws.Cells.Item(44,6).Value = "Active"
